# Apply "Append: 2025-12-13 18:23 JST" update to the ランサーズ sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# Row 2: only the fetch timestamp changes.
$ws.Range("A2").Value = "2025-12-13 18:23:48"

# Row 3: becomes the newly scraped job entry.
$ws.Range("A3").Value = "2025-12-13 18:23:48"
$ws.Range("B3").Value = "【Goエンジニア】OSS活動支援の依頼"
$ws.Range("D3").Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5453259"
$ws.Range("G3").Value = 10
